# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 1
    4  = 3
    5  = 2
    6  = 2
    7  = 1
    8  = 2
    9  = 1
    10 = 1
    11 = 0
    12 = 1
    13 = 3
    14 = 2
    15 = 1
    16 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
